$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Rv2247"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = "accD6 Rv2247 MTCY427.28"
$ws.Range("D4").Value = "FUNCTION: Component of a biotin-dependent acyl-CoA carboxylase complex. This subunit transfers the CO2 from carboxybiotin to the CoA ester substrate (PubMed:17114269). When associated with the alpha3 subunit AccA3, is involved in the carboxylation of acetyl-CoA and propionyl-CoA, with a preference for acetyl-CoA (PubMed:17114269). {ECO:0000269|PubMed:17114269}."
$ws.Range("E4").Value = 54

$ws.Range("A5").Value = "Rv2921c"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "ftsY Rv2921c MTCY338.10c"
$ws.Range("D5").Value = "FUNCTION: Involved in targeting and insertion of nascent membrane proteins into the cytoplasmic membrane. Acts as a receptor for the complex formed by the signal recognition particle (SRP) and the ribosome-nascent chain (RNC) (By similarity). Most of the substrate proteins are involved in stress regulation, lipid metabolism, intermediary metabolism, and cell wall processes (PubMed:29361248). Shows GTPase activity (PubMed:29361248, PubMed:33412199). Can also hydrolyze ATP, UTP and CTP (PubMed:33412199). {ECO:0000255|HAMAP-Rule:MF_00920, ECO:0000269|PubMed:29361248, ECO:0000269|PubMed:33412199}."
$ws.Range("E5").Value = 54

$ws.Range("A6").Value = "Rv2676c"
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "chdC hemQ Rv2676c"
$ws.Range("D6").Value = "FUNCTION: Involved in coproporphyrin-dependent heme b biosynthesis (PubMed:25646457). Catalyzes the decarboxylation of Fe-coproporphyrin III (coproheme) to heme b (protoheme IX), the last step of the pathway (PubMed:25646457). The reaction occurs in a stepwise manner with a three-propionate intermediate (By similarity). {ECO:0000250|UniProtKB:Q8Y5F1, ECO:0000269|PubMed:25646457}."
$ws.Range("E6").Value = 54

$ws.Range("A7").Value = "Rv1689"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "tyrS Rv1689 MTCI125.11"
$ws.Range("D7").Value = "FUNCTION: Catalyzes the attachment of tyrosine to tRNA(Tyr) in a two-step reaction: tyrosine is first activated by ATP to form Tyr-AMP and then transferred to the acceptor end of tRNA(Tyr). {ECO:0000255|HAMAP-Rule:MF_02006}."
$ws.Range("E7").Value = 54

$ws.Range("A8").Value = "Rv1547"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "dnaE1 dnaE Rv1547 MTCY48.18c"
$ws.Range("D8").Value = "FUNCTION: DNA polymerase III is a complex, multichain enzyme responsible for most of the replicative synthesis in bacteria. This DNA polymerase also exhibits 3' to 5' exonuclease activity. The alpha chain is the DNA polymerase (By similarity). {ECO:0000250}."
$ws.Range("E8").Value = 54

$ws.Range("A9").Value = "Rv0684"
$ws.Range("B9").Value = 3
$ws.Range("C9").Value = "fusA Rv0684 MTCY210.01 MTV040.12"
$ws.Range("D9").Value = "FUNCTION: Catalyzes the GTP-dependent ribosomal translocation step during translation elongation. During this step, the ribosome changes from the pre-translocational (PRE) to the post-translocational (POST) state as the newly formed A-site-bound peptidyl-tRNA and P-site-bound deacylated tRNA move to the P and E sites, respectively. Catalyzes the coordinated movement of the two tRNA molecules, the mRNA and conformational changes in the ribosome (By similarity). {ECO:0000250}."
$ws.Range("E9").Value = 54

$ws.Range("A10").Value = "Rv1990c"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = "mbcA Rv1990c MTCY39.29"
$ws.Range("D10").Value = "FUNCTION: Antitoxin component of a type II toxin-antitoxin (TA) system (PubMed:30315706, PubMed:30792174). Neutralizes the activity of cognate toxin MbcT by blocking access to the toxin active site (PubMed:30792174). {ECO:0000269|PubMed:30315706, ECO:0000269|PubMed:30792174}."
$ws.Range("E10").Value = 54

$ws.Range("A11").Value = "Rv0718"
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = "rpsH Rv0718 MTCY210.37"
$ws.Range("D11").Value = "FUNCTION: One of the primary rRNA binding proteins, it binds directly to 16S rRNA central domain where it helps coordinate assembly of the platform of the 30S subunit. {ECO:0000255|HAMAP-Rule:MF_01302}."
$ws.Range("E11").Value = 54

$ws.Range("A12").Value = "Rv0528"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Rv0528"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = 54
